# Update the GitHub->GitLab migration sheet: source is now GitHub ("code-migration"
# username / "repo-migration" target namespace) with a full list of CASA repos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing rows: refresh values for the new source/target pair ---
$ws.Range("B2").Value = "code-migration"
$ws.Range("C2").Value = "casaplotms "
$ws.Range("D2").Value = "repo-migration"

$ws.Range("B3").Value = "code-migration"
$ws.Range("C3").Value = "casa-build-utils "
$ws.Range("D3").Value = "repo-migration"

# --- new rows for the rest of the repo list ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "code-migration"
$ws.Range("C4").Value = "casashell"
$ws.Range("D4").Value = "repo-migration"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "code-migration"
$ws.Range("C5").Value = "casaaddons "
$ws.Range("D5").Value = "repo-migration"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "code-migration"
$ws.Range("C6").Value = "cartavis "
$ws.Range("D6").Value = "repo-migration"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "code-migration"
$ws.Range("C7").Value = "carta-casacore"
$ws.Range("D7").Value = "repo-migration"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "code-migration"
$ws.Range("C8").Value = "casa-asap "
$ws.Range("D8").Value = "repo-migration"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "code-migration"
$ws.Range("C9").Value = "almatasks "
$ws.Range("D9").Value = "repo-migration"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "code-migration"
$ws.Range("C10").Value = "app-n-pak "
$ws.Range("D10").Value = "repo-migration"

# C4 ("casashell") picked up the built-in Hyperlink cell style (left/top aligned,
# vertically centered, wrapped) while the edit was made.
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C4").HorizontalAlignment = -4131
$ws.Range("C4").VerticalAlignment = -4108
$ws.Range("C4").WrapText = $true

# Leave the selection where the editor left off.
$ws.Range("D4").Select() | Out-Null
